# ProductOrders.xlsx edit: remove stale "Failed order" demo row from Orders,
# and add a new Invoice sheet (UnitPrice / Total columns) used to prototype
# the Project3 email-invoice logic.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Orders sheet: drop the D column ("Notes") entirely, and drop the
#    failed-order demo row's Status/Notes cells.
# ---------------------------------------------------------------------
$orders = $wb.Worksheets.Item("Orders")
$orders.Range("D2:D9").ClearContents()
$orders.Range("C9:D9").ClearContents()
$orders.Range("C25").Select()

# ---------------------------------------------------------------------
# 2) New "Invoice" sheet, appended after "Address".
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$inv = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$inv.Name = "Invoice"

$inv.Columns.Item(1).ColumnWidth = 26.75
$inv.Columns.Item(2).ColumnWidth = 9.3
$inv.Columns.Item(4).ColumnWidth = 12.1

$currencyFormat = '"$"#,##0.00'

# Headers
$inv.Range("A1").Value = "Product"
$inv.Range("C1").Value = "Quantity"

$inv.Range("B1:B8").NumberFormat = $currencyFormat
$inv.Range("D1:D8").NumberFormat = $currencyFormat

$inv.Range("B1").Value = "UnitPrice"
$inv.Range("D1").Value = "Total"

# Product names (column A) and quantities (column C)
$inv.Range("A2").Value = "Chai"
$inv.Range("A3").Value = "Ipoh Coffee"
$inv.Range("A4").Value = "Sasquatch Ale"
$inv.Range("A5").Value = "Outback Lager"
$inv.Range("A6").Value = "Guarana Fantastica"
$inv.Range("A7").Value = "Steeleye Stout"
$inv.Range("A8").Value = "Laughing Lumberjack Lager"

$inv.Range("C2").Value = 2
$inv.Range("C3").Value = 3
$inv.Range("C4").Value = 4
$inv.Range("C5").Value = 1
$inv.Range("C6").Value = 2
$inv.Range("C7").Value = 10
$inv.Range("C8").Value = 2

# UnitPrice / Total values need to land as literal currency-formatted TEXT
# (shared strings), matching how they were pasted in from a TEXT()-formula
# staging area rather than typed as numbers. Use a helper cell off to the
# side: build the text with TEXT(), copy it, and paste-special the VALUE
# only into the destination (which already carries the currency style) so
# the destination's number format/style is left untouched.
function Set-InvoiceText($cell, $amount) {
    $helper = $inv.Range("Z1")
    $helper.Formula = '=TEXT(' + $amount + ',"$0.00")'
    $helper.Copy()
    $inv.Range($cell).PasteSpecial(-4163)
}

# Fill order matches the real authoring order: headers, then UnitPrice for
# rows 2-7, then Total for rows 2-7, then the last row (8) UnitPrice/Total.
Set-InvoiceText "B2" 10
Set-InvoiceText "B3" 23
Set-InvoiceText "B4" 7
Set-InvoiceText "B5" 7.5
Set-InvoiceText "B6" 2.25
Set-InvoiceText "B7" 9

Set-InvoiceText "D2" 20
Set-InvoiceText "D3" 69
Set-InvoiceText "D4" 28
Set-InvoiceText "D5" 7.5
Set-InvoiceText "D6" 4.5
Set-InvoiceText "D7" 90

Set-InvoiceText "B8" 7
Set-InvoiceText "D8" 14

$inv.Range("Z1").Clear()

$inv.PageSetup.Orientation = 1

$inv.Range("D13").Select()
$inv.Activate()
